# Update the "Eval_Operador" (T) column scores on Hoja1 to reflect the
# newly-added evaluation criterion (column V), add the new column's width,
# and move the active selection to M7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Updated Eval_Operador values (column T) ---
$ws.Range("T2").Value  = 3.8571428571428572
$ws.Range("T4").Value  = 4.7714285714285714
$ws.Range("T7").Value  = 2.7999999999999994
$ws.Range("T8").Value  = 4.7714285714285714
$ws.Range("T10").Value = 4.628571428571429
$ws.Range("T11").Value = 4.7714285714285714
$ws.Range("T12").Value = 4.5999999999999996
$ws.Range("T13").Value = 3.7142857142857144
$ws.Range("T14").Value = 2.7999999999999994
$ws.Range("T17").Value = 2.7999999999999994
$ws.Range("T18").Value = 3.6857142857142855
$ws.Range("T21").Value = 4.5
$ws.Range("T22").Value = 2.7999999999999994
$ws.Range("T23").Value = 2.7999999999999994
$ws.Range("T25").Value = 4.0285714285714285
$ws.Range("T26").Value = 3.7428571428571429
$ws.Range("T29").Value = 4.5
$ws.Range("T31").Value = 4.7714285714285714
$ws.Range("T32").Value = 4.5
$ws.Range("T33").Value = 4.628571428571429
$ws.Range("T35").Value = 2.7999999999999994
$ws.Range("T36").Value = 4.5
$ws.Range("T38").Value = 3.8571428571428572
$ws.Range("T40").Value = 3.8571428571428572
$ws.Range("T42").Value = 4.5
$ws.Range("T43").Value = 4.5
$ws.Range("T44").Value = 3.7428571428571429
$ws.Range("T45").Value = 4.5
$ws.Range("T46").Value = 4.4119047619047622

# --- New column (V) introduced for the extra criterion: give it an explicit
#     width of 12 characters (ColumnWidth is offset from the stored OOXML
#     width by the default ~0.8333 padding, so back that out here). ---
$ws.Columns.Item(22).ColumnWidth = 11.166666666666666

# --- Move / record the active selection on M7 ---
$ws.Range("M7").Select()
